$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 / Row 2 headers stay the same text, just re-assert (values unchanged)
$ws.Range("B1").Value = "# 두번째 줄까지 수정 금지!, 그 외에는 자유롭게 수정하세요."

$ws.Range("B2").Value = "상담사"
$ws.Range("C2").Value = "접수경로"
$ws.Range("E2").Value = "구글 시트 이름"
$ws.Range("F2").Value = "Main 시트"
$ws.Range("G2").Value = "제품 정보 시트"
$ws.Range("H2").Value = "사은품 정보 시트"

# Row 4: consultant name + unchanged 접수경로 value + renamed setting labels
$ws.Range("B4").Value = "consultant1"
$ws.Range("C4").Value = "옥션"
$ws.Range("E4").Value = "google sheet"
$ws.Range("F4").Value = "Main sheet name"
$ws.Range("G4").Value = "model information name"
$ws.Range("H4").Value = "gift information name"

# Row 5/6: remaining consultant names, 접수경로 options unchanged
$ws.Range("B5").Value = "consultant2"
$ws.Range("C5").Value = "G마켓"

$ws.Range("B6").Value = "consultant3"
$ws.Range("C6").Value = "11번가"

# Rows 7-13 in column C: remove the extra 접수경로 options (인터파크, CJ몰, NS몰,
# 롯데 닷컴, 롯데 아이몰, GS, 홈앤쇼핑) - clear the cell contents
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()

# Restore view state: scroll so column B is the left-most visible column,
# set zoom to 100%, and move the active selection to F9.
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("F9").Select()
